# "updated 2 icdc scripts to resolve wait time issue"
#
# The CasesTab query (cell B2 on the single "startup" sheet) used to look up a
# cohort and return an extra `Cohort` column. That cohort lookup/column is
# removed from the RETURN clause (it was adding to the query's wait time).
# The SamplesTab (B3) / FilesTab (B4) query text is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b2 = $ws.Range("B2")
$current = $b2.Value2

# Drop the trailing ", coalesce(co.cohort_description, '') AS `Cohort`" return
# column (and the now-dangling comma before it) from the end of the query.
$cohortSuffix = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"

if ($current.EndsWith($cohortSuffix)) {
    $newText = $current.Substring(0, $current.Length - $cohortSuffix.Length)
    $b2.Value2 = $newText
}

# The query text is now shorter, so the wrapped row shrinks from 319 to 304.5pt.
$ws.Rows.Item(2).RowHeight = 304.5

# Leave the selection / active cell on B2 (top-left scrolled to row 2), as the
# author did after making this edit.
$b2.Select()
